# Apply updated admin-data "average" (column C, row 2 only) and "p"
# (column D, rows 2-9) values for each year sheet (2008-2011) in the ECU
# total-pos workbook. Columns A/B/E/F keep their existing displayed
# values (the underlying shared-string indices shift in the source diff,
# but they resolve to identical text, so no visible change is needed
# there).

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "2008"; Cell = "C2"; Value = 27.642103216567037 },
    @{ Sheet = "2008"; Cell = "D2"; Value = 0.99942106008529663 },
    @{ Sheet = "2008"; Cell = "D3"; Value = 0.99971026182174683 },
    @{ Sheet = "2008"; Cell = "D4"; Value = 0.99994170665740967 },
    @{ Sheet = "2008"; Cell = "D5"; Value = 0.99997061491012573 },
    @{ Sheet = "2008"; Cell = "D6"; Value = 0.99999374151229858 },
    @{ Sheet = "2008"; Cell = "D7"; Value = 0.99999666213989258 },
    @{ Sheet = "2008"; Cell = "D8"; Value = 0.99999898672103882 },
    @{ Sheet = "2008"; Cell = "D9"; Value = 0.99999946355819702 },

    @{ Sheet = "2009"; Cell = "C2"; Value = 27.123078602786062 },
    @{ Sheet = "2009"; Cell = "D2"; Value = 0.99941641092300415 },
    @{ Sheet = "2009"; Cell = "D3"; Value = 0.99970799684524536 },
    @{ Sheet = "2009"; Cell = "D4"; Value = 0.99994122982025146 },
    @{ Sheet = "2009"; Cell = "D5"; Value = 0.99997037649154663 },
    @{ Sheet = "2009"; Cell = "D6"; Value = 0.99999368190765381 },
    @{ Sheet = "2009"; Cell = "D7"; Value = 0.9999966025352478 },
    @{ Sheet = "2009"; Cell = "D8"; Value = 0.99999892711639404 },
    @{ Sheet = "2009"; Cell = "D9"; Value = 0.99999946355819702 },

    @{ Sheet = "2010"; Cell = "C2"; Value = 29.466916212687494 },
    @{ Sheet = "2010"; Cell = "D2"; Value = 0.99941170215606689 },
    @{ Sheet = "2010"; Cell = "D3"; Value = 0.99970561265945435 },
    @{ Sheet = "2010"; Cell = "D4"; Value = 0.99994075298309326 },
    @{ Sheet = "2010"; Cell = "D5"; Value = 0.99997013807296753 },
    @{ Sheet = "2010"; Cell = "D6"; Value = 0.99999362230300903 },
    @{ Sheet = "2010"; Cell = "D7"; Value = 0.9999966025352478 },
    @{ Sheet = "2010"; Cell = "D8"; Value = 0.99999892711639404 },
    @{ Sheet = "2010"; Cell = "D9"; Value = 0.99999946355819702 },

    @{ Sheet = "2011"; Cell = "C2"; Value = 33.465564727783203 },
    @{ Sheet = "2011"; Cell = "D2"; Value = 0.99940687417984009 },
    @{ Sheet = "2011"; Cell = "D3"; Value = 0.99970322847366333 },
    @{ Sheet = "2011"; Cell = "D4"; Value = 0.99994027614593506 },
    @{ Sheet = "2011"; Cell = "D5"; Value = 0.99996989965438843 },
    @{ Sheet = "2011"; Cell = "D6"; Value = 0.99999362230300903 },
    @{ Sheet = "2011"; Cell = "D7"; Value = 0.99999654293060303 },
    @{ Sheet = "2011"; Cell = "D8"; Value = 0.99999892711639404 },
    @{ Sheet = "2011"; Cell = "D9"; Value = 0.99999946355819702 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item([string]$u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
